$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fertig/% (F9) and Zeit/min tatsächlich (G9) for row 9
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 50

# Move active selection from G8 to D8
$ws.Range("D8").Select()
